# Add a new "2023" column (L) to the table, mirroring the formatting of
# column K (2022) on each row, and refresh a couple of row heights that
# Excel recomputed as part of the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column L: header year + 10 data rows -------------------------
# Row 4 holds the year headers (2015 ... 2022); append 2023.
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("L4").Value = 2023

# Data rows 5-14: copy the number format/style from column K, then set
# the 2023 value.
$values = @{
    5  = 1.6430457248453274
    6  = 0.41181606829870221
    7  = 0.94796963217320562
    8  = 0.72306112208737106
    9  = 2.1802539701246277
    10 = 0.63651150401750112
    11 = 0.97994201681774651
    12 = 2.2469385026996971
    13 = 4.1686356866605365
    14 = 0.3304193846038968
}

foreach ($row in 5..14) {
    $srcCell = "K$row"
    $dstCell = "L$row"
    $ws.Range($srcCell).Copy()
    $ws.Range($dstCell).PasteSpecial(-4122)   # xlPasteFormats
    $ws.Range($dstCell).Value = $values[$row]
}

# --- Row height refresh -------------------------------------------------
# Rows 2 and 3 were re-measured (12.75 -> 13.5) as part of the edit.
$ws.Rows.Item(2).RowHeight = 13.5
$ws.Rows.Item(3).RowHeight = 13.5

# --- Reset selection back to the top of the sheet -----------------------
[void]$ws.Range("A1").Select()

Write-Output "done"
